$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.325.73"
$ws.Range("E2").Value = "  +0.07%  "

# Row 3
$ws.Range("D3").Value = "1.931.53"
$ws.Range("E3").Value = "  +0.16%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7463"
$ws.Range("E5").Value = "  +4.24%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.12"
$ws.Range("E6").Value = "  -2.47%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "27.58"
$ws.Range("E8").Value = "  -0.55%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3166"
$ws.Range("E9").Value = "  -1.28%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07138"
$ws.Range("E10").Value = "  +0.52%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08057"
$ws.Range("E11").Value = "  +0.83%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7788"
$ws.Range("E12").Value = "  -1.32%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.404"
$ws.Range("E13").Value = "  +0.31%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.884.10"
$ws.Range("E14").Value = "  -2.33%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.17"
$ws.Range("E15").Value = "  -1.72%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.59"
$ws.Range("E16").Value = "  -0.64%  "

# Row 17
$ws.Range("D17").Value = "30.310.76"
$ws.Range("E17").Value = "  +0.08%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.017"
$ws.Range("E18").Value = "  +4.45%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "251.92"
$ws.Range("E19").Value = "  -2.37%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007906"
$ws.Range("E20").Value = "  -2.36%  "

# Row 21
$ws.Range("E21").Value = "  +0.01%  "

# Row 22
$ws.Range("D22").Value = "2.153.61"
$ws.Range("E22").Value = "  -1.32%  "

# Row 23
$ws.Range("E23").Value = "  +0.07%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.659"
$ws.Range("E24").Value = "  -2.59%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.570"
$ws.Range("E25").Value = "  +0.39%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.57"
$ws.Range("E26").Value = "  +0.39%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.07"
$ws.Range("E27").Value = "  -0.13%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1293"
$ws.Range("E28").Value = "  +2.35%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.178"
$ws.Range("E29").Value = "  -4.11%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.367"
$ws.Range("E30").Value = "  +0.94%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.555"
$ws.Range("E31").Value = "  +1.52%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.417"
$ws.Range("E32").Value = "  +0.36%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.147"
$ws.Range("E33").Value = "  +0.12%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05248"
$ws.Range("E34").Value = "  +1.98%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.317"
$ws.Range("E35").Value = "  +4.01%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7580"
$ws.Range("E36").Value = "  +1.88%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.785"
$ws.Range("E37").Value = "  +0.61%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01953"
$ws.Range("E38").Value = "  -0.77%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.800"
$ws.Range("E39").Value = "  +0.06%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.516"
$ws.Range("E40").Value = "  +2.27%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "77.80"
$ws.Range("E41").Value = "  -0.30%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4529"
$ws.Range("E42").Value = "  +0.41%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.974"
$ws.Range("E43").Value = "  -1.03%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8424"
$ws.Range("E44").Value = "  -0.67%  "

# Row 45
$ws.Range("E45").Value = "  +0.07%  "

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.03"
$ws.Range("E46").Value = "  +1.84%  "

# Row 47
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.701"
$ws.Range("E47").Value = "  +3.47%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.57"
$ws.Range("E48").Value = "  +1.06%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.82"
$ws.Range("E49").Value = "  +2.72%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1219"
$ws.Range("E50").Value = "  +6.90%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "951.17"
$ws.Range("E51").Value = "  -0.02%  "
